$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column holds price values that look numeric (e.g. "0.9991", "289.85").
# Excel would silently coerce these strings to numbers on assignment, which
# would not match the source workbook where they are stored as literal text.
# Temporarily force the Text number format on D2:D51 before writing the new
# values, then restore the default "Normal" style so no stray number-format
# / quote-prefix styling is left behind on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "21.653.56"
$ws.Range("E2").Value = "  -1.82%  "

# Row 3
$ws.Range("D3").Value = "1.535.31"
$ws.Range("E3").Value = "  -1.26%  "

# Row 4
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "0.9998"
$ws.Range("E5").Value = "  -0.06%  "

# Row 6
$ws.Range("D6").Value = "289.85"
$ws.Range("E6").Value = "  +0.66%  "

# Row 7
$ws.Range("E7").Value = "  -0.89%  "

# Row 8
$ws.Range("D8").Value = "0.3171"
$ws.Range("E8").Value = "  -1.86%  "

# Row 9
$ws.Range("D9").Value = "42.33"
$ws.Range("E9").Value = "  -0.90%  "

# Row 10
$ws.Range("D10").Value = "0.07168"
$ws.Range("E10").Value = "  -2.37%  "

# Row 11
$ws.Range("D11").Value = "1.063"
$ws.Range("E11").Value = "  -3.56%  "

# Row 12
$ws.Range("E12").Value = "  +0.09%  "

# Row 13
$ws.Range("D13").Value = "5.708"
$ws.Range("E13").Value = "  +0.40%  "

# Row 14
$ws.Range("D14").Value = "18.33"
$ws.Range("E14").Value = "  -3.30%  "

# Row 15
$ws.Range("D15").Value = "6.624"
$ws.Range("E15").Value = "  -1.59%  "

# Row 16
$ws.Range("D16").Value = "1.537.55"
$ws.Range("E16").Value = "  -0.45%  "

# Row 17
$ws.Range("D17").Value = "0.00001094"
$ws.Range("E17").Value = "  -3.51%  "

# Row 18
$ws.Range("D18").Value = "0.06597"
$ws.Range("E18").Value = "  -0.43%  "

# Row 19
$ws.Range("D19").Value = "84.00"
$ws.Range("E19").Value = "  -1.39%  "

# Row 20
$ws.Range("D20").Value = "0.9995"
$ws.Range("E20").Value = "  -0.06%  "

# Row 21
$ws.Range("D21").Value = "6.148"
$ws.Range("E21").Value = "  -3.05%  "

# Row 22
$ws.Range("D22").Value = "15.55"
$ws.Range("E22").Value = "  -2.26%  "

# Row 23
$ws.Range("D23").Value = "10.73"
$ws.Range("E23").Value = "  -4.84%  "

# Row 24
$ws.Range("D24").Value = "2.370"
$ws.Range("E24").Value = "  +1.50%  "

# Row 25
$ws.Range("D25").Value = "21.618.06"
$ws.Range("E25").Value = "  -1.96%  "

# Row 26
$ws.Range("D26").Value = "2.354"
$ws.Range("E26").Value = "  -4.46%  "

# Row 27
$ws.Range("D27").Value = "150.38"
$ws.Range("E27").Value = "  +1.56%  "

# Row 28
$ws.Range("D28").Value = "18.36"
$ws.Range("E28").Value = "  -2.20%  "

# Row 29
$ws.Range("D29").Value = "4.854"
$ws.Range("E29").Value = "  -0.17%  "

# Row 30
$ws.Range("D30").Value = "1.707.80"
$ws.Range("E30").Value = "  -0.81%  "

# Row 31
$ws.Range("D31").Value = "117.39"
$ws.Range("E31").Value = "  -2.96%  "

# Row 32
$ws.Range("D32").Value = "6.081"
$ws.Range("E32").Value = "  +5.94%  "

# Row 33
$ws.Range("D33").Value = "0.9463"
$ws.Range("E33").Value = "  -11.19%  "

# Row 34
$ws.Range("D34").Value = "0.08155"
$ws.Range("E34").Value = "  -2.66%  "

# Row 35
$ws.Range("D35").Value = "5.180"
$ws.Range("E35").Value = "  +0.74%  "

# Row 36
$ws.Range("D36").Value = "8.492"
$ws.Range("E36").Value = "  -8.82%  "

# Row 37
$ws.Range("D37").Value = "0.06006"
$ws.Range("E37").Value = "  -3.70%  "

# Row 38
$ws.Range("D38").Value = "0.02220"
$ws.Range("E38").Value = "  -2.58%  "

# Row 39
$ws.Range("D39").Value = "1.461"
$ws.Range("E39").Value = "  -10.60%  "

# Row 40
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "0.2035"
$ws.Range("E40").Value = "  -3.17%  "

# Row 41
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "11.09"
$ws.Range("E41").Value = "  +3.64%  "

# Row 42
$ws.Range("D42").Value = "1.185"
$ws.Range("E42").Value = "  -2.30%  "

# Row 43
$ws.Range("D43").Value = "0.9990"
$ws.Range("E43").Value = "  -0.11%  "

# Row 44
$ws.Range("D44").Value = "0.5810"
$ws.Range("E44").Value = "  -1.30%  "

# Row 45
$ws.Range("D45").Value = "12.99"
$ws.Range("E45").Value = "  -2.05%  "

# Row 46
$ws.Range("D46").Value = "3.718"
$ws.Range("E46").Value = "  +0.01%  "

# Row 47
$ws.Range("D47").Value = "0.5557"
$ws.Range("E47").Value = "  -1.09%  "

# Row 48
$ws.Range("E48").Value = "  +1.62%  "

# Row 49
$ws.Range("D49").Value = "1.881"
$ws.Range("E49").Value = "  -1.81%  "

# Row 50
$ws.Range("D50").Value = "116.40"
$ws.Range("E50").Value = "  -1.46%  "

# Row 51
$ws.Range("D51").Value = "0.06699"
$ws.Range("E51").Value = "  -2.32%  "

# Restore default styling on the price column now that the text values are set
$ws.Range("D2:D51").Style = "Normal"